$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 849  # F2: 848 -> 849
$ws.Cells.Item(3, 6).Value = 1755  # F3: 1750 -> 1755
$ws.Cells.Item(4, 6).Value = 46  # F4: 43 -> 46
$ws.Cells.Item(5, 6).Value = 533  # F5: 532 -> 533
$ws.Cells.Item(6, 6).Value = 2145  # F6: 2146 -> 2145
$ws.Cells.Item(7, 6).Value = 1367  # F7: 1365 -> 1367
$ws.Cells.Item(8, 6).Value = 2073  # F8: 2065 -> 2073
$ws.Cells.Item(11, 6).Value = 2403  # F11: 2401 -> 2403
$ws.Cells.Item(12, 6).Value = 658  # F12: 657 -> 658
$ws.Cells.Item(14, 6).Value = 3914  # F14: 3907 -> 3914
$ws.Cells.Item(17, 6).Value = 2992  # F17: 2982 -> 2992
$ws.Cells.Item(18, 6).Value = 800  # F18: 791 -> 800
$ws.Cells.Item(19, 6).Value = 142  # F19: 140 -> 142
$ws.Cells.Item(20, 6).Value = 1345  # F20: 1344 -> 1345
$ws.Cells.Item(21, 6).Value = 112  # F21: 108 -> 112
$ws.Cells.Item(22, 6).Value = 2047  # F22: 2044 -> 2047
$ws.Cells.Item(23, 6).Value = 1176  # F23: 1171 -> 1176
$ws.Cells.Item(24, 6).Value = 1865  # F24: 1850 -> 1865
$ws.Cells.Item(25, 6).Value = 383  # F25: 382 -> 383
$ws.Cells.Item(26, 6).Value = 206  # F26: 204 -> 206
$ws.Cells.Item(27, 6).Value = 11  # F27: 10 -> 11
$ws.Cells.Item(28, 6).Value = 8361  # F28: 8322 -> 8361
$ws.Cells.Item(29, 6).Value = 5647  # F29: 5627 -> 5647
$ws.Cells.Item(30, 6).Value = 353  # F30: 352 -> 353
$ws.Cells.Item(31, 6).Value = 176  # F31: 172 -> 176
$ws.Cells.Item(33, 6).Value = 765  # F33: 762 -> 765
$ws.Cells.Item(36, 6).Value = 945  # F36: 942 -> 945
$ws.Cells.Item(37, 6).Value = 387  # F37: 381 -> 387
$ws.Cells.Item(38, 6).Value = 34  # F38: 33 -> 34
$ws.Cells.Item(39, 6).Value = 189  # F39: 188 -> 189
$ws.Cells.Item(40, 6).Value = 159  # F40: 154 -> 159
$ws.Cells.Item(41, 6).Value = 4617  # F41: 4613 -> 4617
$ws.Cells.Item(42, 6).Value = 848  # F42: 841 -> 848
$ws.Cells.Item(43, 6).Value = 73  # F43: 71 -> 73

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 87  # F3: 86 -> 87
$ws.Cells.Item(16, 6).Value = 25  # F16: 24 -> 25
$ws.Cells.Item(18, 6).Value = 170  # F18: 169 -> 170

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 8264  # F2: 8246 -> 8264
$ws.Cells.Item(3, 6).Value = 361  # F3: 358 -> 361
$ws.Cells.Item(4, 6).Value = 1271  # F4: 1264 -> 1271

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 849  # F2: 848 -> 849
$ws.Cells.Item(3, 6).Value = 361  # F3: 358 -> 361
$ws.Cells.Item(4, 6).Value = 1271  # F4: 1264 -> 1271
$ws.Cells.Item(5, 6).Value = 87  # F5: 86 -> 87
$ws.Cells.Item(6, 6).Value = 1755  # F6: 1750 -> 1755
$ws.Cells.Item(7, 6).Value = 46  # F7: 43 -> 46
$ws.Cells.Item(8, 6).Value = 533  # F8: 532 -> 533
$ws.Cells.Item(9, 6).Value = 1367  # F9: 1365 -> 1367
$ws.Cells.Item(10, 6).Value = 2073  # F10: 2065 -> 2073
$ws.Cells.Item(15, 6).Value = 3914  # F15: 3907 -> 3914
$ws.Cells.Item(17, 6).Value = 2992  # F17: 2982 -> 2992
$ws.Cells.Item(18, 6).Value = 800  # F18: 791 -> 800
$ws.Cells.Item(19, 6).Value = 142  # F19: 140 -> 142
$ws.Cells.Item(21, 6).Value = 2047  # F21: 2044 -> 2047
$ws.Cells.Item(27, 6).Value = 1865  # F27: 1850 -> 1865
$ws.Cells.Item(29, 6).Value = 206  # F29: 204 -> 206
$ws.Cells.Item(30, 6).Value = 11  # F30: 10 -> 11
$ws.Cells.Item(31, 6).Value = 8361  # F31: 8322 -> 8361
$ws.Cells.Item(32, 6).Value = 5647  # F32: 5627 -> 5647
$ws.Cells.Item(34, 6).Value = 353  # F34: 352 -> 353
$ws.Cells.Item(35, 6).Value = 176  # F35: 172 -> 176
$ws.Cells.Item(37, 6).Value = 765  # F37: 762 -> 765
$ws.Cells.Item(39, 6).Value = 945  # F39: 942 -> 945
$ws.Cells.Item(40, 6).Value = 387  # F40: 381 -> 387
$ws.Cells.Item(41, 6).Value = 189  # F41: 188 -> 189
$ws.Cells.Item(42, 6).Value = 159  # F42: 154 -> 159
$ws.Cells.Item(43, 6).Value = 4617  # F43: 4613 -> 4617
$ws.Cells.Item(44, 6).Value = 848  # F44: 841 -> 848
